$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the product list with additional empty rows (16-25), continuing the
# sequential Id numbering in column A (15..24) and leaving the rest of each
# row's columns (Categoria, Comida, Cantidad Disponible, Precio) blank text
# cells, matching the existing "placeholder" rows already present (11-15).
for ($i = 16; $i -le 25; $i++) {
    $ws.Cells.Item($i, 1).Value = $i - 1
    for ($c = 2; $c -le 5; $c++) {
        $cell = $ws.Cells.Item($i, $c)
        # Writing a bare "'" materializes an empty *text* cell (Excel's
        # quote-prefix literal-text marker) instead of leaving the cell
        # completely absent/blank. Resetting the style back to Normal
        # afterwards drops the quote-prefix formatting, leaving a plain
        # empty inline/shared string cell like the pre-existing rows.
        $cell.Value = "'"
        $cell.Style = "Normal"
    }
}
